$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the confidential disclosure date in A38
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for holdings rows 2-35
$ws.Range("D2").Value = 0.03550564198204227
$ws.Range("E2").Value = 0
$ws.Range("D3").Value = 0.02017317617660333
$ws.Range("E3").Value = -0.001970055161544582
$ws.Range("D4").Value = 0.01919123640150408
$ws.Range("E4").Value = 0.001779661016949152
$ws.Range("D5").Value = 0.03725175427410597
$ws.Range("E5").Value = 0.009601706970128188
$ws.Range("D6").Value = 0.03425126385194842
$ws.Range("E6").Value = 0
$ws.Range("D7").Value = 0.01977800808762834
$ws.Range("E7").Value = -0.0008945770853017176
$ws.Range("D8").Value = 0.03740257109302948
$ws.Range("E8").Value = -0.004856418918918859
$ws.Range("D9").Value = 0.02040245850266948
$ws.Range("E9").Value = -0.0005394228175852334
$ws.Range("D10").Value = 0.02625537391597386
$ws.Range("E10").Value = -0.005918882204541065
$ws.Range("D11").Value = 0.02410648900445397
$ws.Range("E11").Value = 0.004491413474240424
$ws.Range("D12").Value = 0.05766771495819935
$ws.Range("E12").Value = -0.001179523472517152
$ws.Range("D13").Value = 0.02464326440558882
$ws.Range("E13").Value = -0.001116486788239834
$ws.Range("D14").Value = 0.02720317610030327
$ws.Range("E14").Value = -0.01029027799109195
$ws.Range("D15").Value = 0.03349193173964538
$ws.Range("E15").Value = -0.01277900834895207
$ws.Range("D16").Value = 0.01977911628553766
$ws.Range("E16").Value = 0.002248875562218755
$ws.Range("D17").Value = 0.0311492777867856
$ws.Range("E17").Value = 0.006788255500122631
$ws.Range("D18").Value = 0.04179205384649572
$ws.Range("E18").Value = 0.0009265693768818473
$ws.Range("D19").Value = 0.1253265346538008
$ws.Range("E19").Value = 0
$ws.Range("D20").Value = 0.009174706801722618
$ws.Range("E20").Value = 0.0201813769319199
$ws.Range("D21").Value = 0.01541633218515113
$ws.Range("E21").Value = -0.008606329816768565
$ws.Range("D22").Value = 0.01714685327904233
$ws.Range("E22").Value = 0.001565972923822878
$ws.Range("D23").Value = 0.01540048622883856
$ws.Range("E23").Value = -0.002884962134872104
$ws.Range("D24").Value = 0.02132422440508337
$ws.Range("E24").Value = 0.0009246892016852026
$ws.Range("D25").Value = 0.01265139113930553
$ws.Range("E25").Value = 0.01268213707501364
$ws.Range("D26").Value = 0.04250435760614119
$ws.Range("E26").Value = -0.00204444923940017
$ws.Range("D27").Value = 0.02389947253983025
$ws.Range("E27").Value = -0.0000980680592330696
$ws.Range("D28").Value = 0.04540768327367471
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0.05584461475620543
$ws.Range("E29").Value = -0.002682403433476366
$ws.Range("D30").Value = 0.01315214373943399
$ws.Range("E30").Value = -0.003858520900321571
$ws.Range("D31").Value = 0.02058878860631922
$ws.Range("E31").Value = -0.0003835826620637306
$ws.Range("D32").Value = 0.01326016119082514
$ws.Range("E32").Value = 0.009125840537944185
$ws.Range("D33").Value = 0.04189813513602231
$ws.Range("E33").Value = -0.001543209876543217
$ws.Range("D34").Value = 0.01695960604608832
$ws.Range("E34").Value = -0.0004416310908287313
$ws.Range("E35").Value = -0.0004479376710774252

$ws.Protect()
